$d = $word.ActiveDocument

function Replace-Text($range, $oldText, $newText) {
    $range.Find.ClearFormatting()
    $ok = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $oldText"
    }
}

# --- 1. "Core Functionality" intro paragraph (para 4): replace text, add spacing ---
$p4 = $d.Paragraphs(4)
Replace-Text $p4.Range "Raspberry pi portable device that connects to a car via OBD2 and have some sort of interface so the user can see stats about their car both real-time and after the journey on a computer for more detailed information" "The core functionality of the device includes:"
$p4 = $d.Paragraphs(4)
$p4.SpaceAfter = 6
$p4.LineSpacingRule = 0
$p4.LineSpacing = 12

# --- 2. Bullet list under "Core Functionality" (paras 5-11): shift content up one slot ---
Replace-Text $d.Paragraphs(5).Range "Connect it to the car, go for a drive get some data together, unplug, connect to network, submit results to web server." "Record journeys in car with GPS;"
Replace-Text $d.Paragraphs(6).Range "Record journeys in car with GPS (attach to raspberry pi usb)" "Dashboard online service (upload to web server)"
Replace-Text $d.Paragraphs(7).Range "Dashboard online service (upload to web server)" "Economic Driving"
Replace-Text $d.Paragraphs(8).Range "Economic Driving" "Can check detailed information like pistons"
Replace-Text $d.Paragraphs(9).Range "Can check detailed information like pistons" "Compare friends results (back end, web server)"

# Delete the old "Compare friends results" bullet paragraph entirely (its text now lives in para 9)
$d.Paragraphs(10).Range.Delete()

# --- 3. "GPS Record Journeys" heading (para 11): Heading1 -> Heading2 ---
$d.Paragraphs(11).Style = "Heading 2"

# --- 4. Para 13: extend text describing colour-coded route ---
Replace-Text $d.Paragraphs(13).Range " and the colour of the route could gradually change depending on speed or how economic you were driving. This is a similar idea to the route colouring technique used in the Nike Running app:" " on a map. The colour of the route will gradually change depending on how economically you were driving. This is a similar idea to the route colouring technique used in the Nike Running app:"

# --- 5. "Dashboard Online Service" heading (para 17): Heading1 -> Heading2 ---
$d.Paragraphs(17).Style = "Heading 2"

# --- 6. Para 18: expand text about the Dashboard Online Service ---
Replace-Text $d.Paragraphs(18).Range "The idea is that when drives have been recorded, our device can be connected to the web service and upload all the information gathered. As a user you will be able to log into the web service and view an analysis of your drive." "The idea is that when drives have been recorded, the device can be connected to the Dashboard Online Service web service and upload all the information gathered. As a user you will be able to log into the web service and view the analysis of your drive or all of your drives."

# --- 7. New paragraph after para 18, describing per-drive statistics ---
$p18 = $d.Paragraphs(18)
$p18.Range.InsertParagraphAfter()
$pNew = $d.Paragraphs(19)
$pNew.Range.Text = "You will be able to view statistics for individual drives and the route for each will be plotted onto a map. There will also be a screen of general statistics and averages of all your drives."
$pNew.Style = "Normal"

# --- 8. "Economic Driving" heading (para 20): Heading1 -> Heading2 ---
$d.Paragraphs(20).Style = "Heading 2"

# --- 9. New paragraph after para 23 ("We will use these factors...") describing the score ---
$p23 = $d.Paragraphs(23)
$p23.Range.InsertParagraphAfter()
$pNew2 = $d.Paragraphs(24)
$pNew2.Range.Text = "An Economic Driving Score will be calculated based on all the drives by a single driver."
$pNew2.Style = "Normal"

# --- 10. "Detailed Information" heading (para 25): Heading1 -> Heading2 ---
$d.Paragraphs(25).Style = "Heading 2"

# --- 11. Empty paragraph after it (para 26) gets placeholder text "?" ---
$p26 = $d.Paragraphs(26)
$p26.Range.Text = "?"
$p26.Style = "Normal"

# --- 12. "Compared Results" heading (para 27): Heading1 -> Heading2 ---
$d.Paragraphs(27).Style = "Heading 2"

# --- 13. Empty paragraph (para 28, carries the _GoBack bookmark) gets the new table description ---
$p28 = $d.Paragraphs(28)
$p28.Range.Text = "There will be a table on the Online Dashboard Service that allows you to compare your general statistics and your economic driving score to others who have used the device. "
$p28.Style = "Normal"

# --- 14. "Real-time information" heading (para 29): Heading1 -> Heading2 ---
$d.Paragraphs(29).Style = "Heading 2"

# Relocate the _GoBack bookmark so it again sits right at the "Real-time information"
# heading / body boundary (it previously lived in the now-repurposed para 28).
$bmRange = $d.Paragraphs(29).Range.Duplicate
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output "$i : [$($d.Paragraphs($i).Range.Text)] style=[$($d.Paragraphs($i).Style.NameLocal)]"
}
